$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 59

# Row 3
$ws.Range("B3").Value = "<was>"
$ws.Range("C3").Value = 59

# Row 4
$ws.Range("B4").Value = "<wor>"
$ws.Range("C4").Value = 59

# Row 5
$ws.Range("B5").Value = "<wand>"
$ws.Range("C5").Value = 60

# Row 6
$ws.Range("B6").Value = "<form>"
$ws.Range("C6").Value = 60

# Row 7
$ws.Range("B7").Value = "<see>"
$ws.Range("C7").Value = 62

# Row 8
$ws.Range("C8").Value = 60

# Row 9
$ws.Range("B9").Value = "<than>"

# Row 10
$ws.Range("B10").Value = "<some>"
$ws.Range("C10").Value = 61

# Row 11
$ws.Range("B11").Value = "<word>"
$ws.Range("C11").Value = 63

# Row 13
$ws.Range("B13").Value = "<his>"
$ws.Range("C13").Value = 61

# Row 14
$ws.Range("B14").Value = "<a>"
